# Auto-génération des classes et des specs
# Simplify the Heading1 titles: drop the "Objet RS-RI:15-15:" prefix and the
# "Type " prefixes, leaving only the bare schema/type name.

$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "Objet RS-RI:15-15:resourcesInfo"; New = "resourcesInfo" },
    @{ Old = "Type resource"; New = "resource" },
    @{ Old = "Type team"; New = "team" },
    @{ Old = "Type state"; New = "state" },
    @{ Old = "Type coord"; New = "coord" },
    @{ Old = "Type contact"; New = "contact" }
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.ClearFormatting()
    $range.Find.Execute($r.Old, $true, $true, $false, $false, $false, $true, 1, $false, $r.New, 2)
}

$d.Save()
